$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cell = $ws.Range("D2")
$cell.NumberFormat = "@"
$cell.Value = "55.384.64"
$cell.Style = "Normal"
$ws.Range("E2").Value = "  -6.17%  "

$cell = $ws.Range("D3")
$cell.NumberFormat = "@"
$cell.Value = "2.923.68"
$cell.Style = "Normal"
$ws.Range("E3").Value = "  -9.67%  "

$cell = $ws.Range("D4")
$cell.NumberFormat = "@"
$cell.Value = "0.999"
$cell.Style = "Normal"
$ws.Range("E4").Value = "  -0.07%  "

$cell = $ws.Range("D5")
$cell.NumberFormat = "@"
$cell.Value = "469.27"
$cell.Style = "Normal"
$ws.Range("E5").Value = "  -12.94%  "

$cell = $ws.Range("D6")
$cell.NumberFormat = "@"
$cell.Value = "123.35"
$cell.Style = "Normal"
$ws.Range("E6").Value = "  -9.76%  "

$ws.Range("E7").Value = "  +0.03%  "

$cell = $ws.Range("D8")
$cell.NumberFormat = "@"
$cell.Value = "2.926.46"
$cell.Style = "Normal"
$ws.Range("E8").Value = "  -9.56%  "

$cell = $ws.Range("D9")
$cell.NumberFormat = "@"
$cell.Value = "0.400"
$cell.Style = "Normal"
$ws.Range("E9").Value = "  -13.01%  "

$cell = $ws.Range("D10")
$cell.NumberFormat = "@"
$cell.Value = "6.57"
$cell.Style = "Normal"
$ws.Range("E10").Value = "  -13.62%  "

$cell = $ws.Range("D11")
$cell.NumberFormat = "@"
$cell.Value = "0.0950"
$cell.Style = "Normal"
$ws.Range("E11").Value = "  -17.75%  "

$cell = $ws.Range("D12")
$cell.NumberFormat = "@"
$cell.Value = "0.328"
$cell.Style = "Normal"
$ws.Range("E12").Value = "  -17.42%  "

$cell = $ws.Range("D13")
$cell.NumberFormat = "@"
$cell.Value = "0.123"
$cell.Style = "Normal"
$ws.Range("E13").Value = "  -3.57%  "

$cell = $ws.Range("D14")
$cell.NumberFormat = "@"
$cell.Value = "3.426.78"
$cell.Style = "Normal"
$ws.Range("E14").Value = "  -9.48%  "

$cell = $ws.Range("D15")
$cell.NumberFormat = "@"
$cell.Value = "22.49"
$cell.Style = "Normal"
$ws.Range("E15").Value = "  -14.06%  "

$cell = $ws.Range("D16")
$cell.NumberFormat = "@"
$cell.Value = "55.516.84"
$cell.Style = "Normal"
$ws.Range("E16").Value = "  -6.10%  "

$cell = $ws.Range("D17")
$cell.NumberFormat = "@"
$cell.Value = "2.925.77"
$cell.Style = "Normal"
$ws.Range("E17").Value = "  -9.39%  "

$ws.Range("E18").Value = "  -17.37%  "

$cell = $ws.Range("D19")
$cell.NumberFormat = "@"
$cell.Value = "5.08"
$cell.Style = "Normal"
$ws.Range("E19").Value = "  -13.77%  "

$cell = $ws.Range("D20")
$cell.NumberFormat = "@"
$cell.Value = "11.44"
$cell.Style = "Normal"
$ws.Range("E20").Value = "  -13.58%  "

$ws.Range("E21").Value = "  -15.72%  "

$cell = $ws.Range("D22")
$cell.NumberFormat = "@"
$cell.Value = "306.10"
$cell.Style = "Normal"
$ws.Range("E22").Value = "  -15.51%  "

$ws.Range("E23").Value = "  -0.29%  "

$cell = $ws.Range("D24")
$cell.NumberFormat = "@"
$cell.Value = "0.443"
$cell.Style = "Normal"
$ws.Range("E24").Value = "  -14.82%  "

$ws.Range("E25").Value = "  -16.79%  "

$ws.Range("E26").Value = "  +0.45%  "

$cell = $ws.Range("D27")
$cell.NumberFormat = "@"
$cell.Value = "0.153"
$cell.Style = "Normal"
$ws.Range("E27").Value = "  -10.13%  "

$cell = $ws.Range("D28")
$cell.NumberFormat = "@"
$cell.Value = "0.997"
$cell.Style = "Normal"
$ws.Range("E28").Value = "  -0.20%  "

$cell = $ws.Range("D29")
$cell.NumberFormat = "@"
$cell.Value = "0.0₃0799"
$cell.Style = "Normal"
$ws.Range("E29").Value = "  -18.33%  "

$cell = $ws.Range("D30")
$cell.NumberFormat = "@"
$cell.Value = "5.88"
$cell.Style = "Normal"
$ws.Range("E30").Value = "  -16.91%  "

$cell = $ws.Range("D31")
$cell.NumberFormat = "@"
$cell.Value = "1.11"
$cell.Style = "Normal"
$ws.Range("E31").Value = "  -9.98%  "

$cell = $ws.Range("D32")
$cell.NumberFormat = "@"
$cell.Value = "19.02"
$cell.Style = "Normal"
$ws.Range("E32").Value = "  -13.90%  "

$cell = $ws.Range("D33")
$cell.NumberFormat = "@"
$cell.Value = "6.01"
$cell.Style = "Normal"
$ws.Range("E33").Value = "  -15.12%  "

$ws.Range("E34").Value = "  -18.79%  "

$cell = $ws.Range("D35")
$cell.NumberFormat = "@"
$cell.Value = "144.51"
$cell.Style = "Normal"
$ws.Range("E35").Value = "  -11.93%  "

$cell = $ws.Range("D36")
$cell.NumberFormat = "@"
$cell.Value = "4.15"
$cell.Style = "Normal"
$ws.Range("E36").Value = "  -15.94%  "

$ws.Range("E37").Value = "  -16.21%  "

$ws.Range("E38").Value = "  -16.17%  "

$cell = $ws.Range("D39")
$cell.NumberFormat = "@"
$cell.Value = "2.957.32"
$cell.Style = "Normal"
$ws.Range("E39").Value = "  -9.42%  "

$cell = $ws.Range("D40")
$cell.NumberFormat = "@"
$cell.Value = "0.999"
$cell.Style = "Normal"
$ws.Range("E40").Value = "  -0.06%  "

$ws.Range("E41").Value = "  -14.73%  "

$cell = $ws.Range("D42")
$cell.NumberFormat = "@"
$cell.Value = "21.49"
$cell.Style = "Normal"
$ws.Range("E42").Value = "  -17.55%  "

$cell = $ws.Range("D43")
$cell.NumberFormat = "@"
$cell.Value = "34.94"
$cell.Style = "Normal"
$ws.Range("E43").Value = "  -15.06%  "

$cell = $ws.Range("D44")
$cell.NumberFormat = "@"
$cell.Value = "0.955"
$cell.Style = "Normal"
$ws.Range("E44").Value = "  -12.71%  "

$ws.Range("E45").Value = "  -16.74%  "

$cell = $ws.Range("D46")
$cell.NumberFormat = "@"
$cell.Value = "3.39"
$cell.Style = "Normal"
$ws.Range("E46").Value = "  -15.86%  "

$cell = $ws.Range("D47")
$cell.NumberFormat = "@"
$cell.Value = "2.048.52"
$cell.Style = "Normal"
$ws.Range("E47").Value = "  -10.79%  "

$ws.Range("E48").Value = "  -14.24%  "

$cell = $ws.Range("D49")
$cell.NumberFormat = "@"
$cell.Value = "5.26"
$cell.Style = "Normal"
$ws.Range("E49").Value = "  -16.40%  "

$cell = $ws.Range("D50")
$cell.NumberFormat = "@"
$cell.Value = "17.46"
$cell.Style = "Normal"
$ws.Range("E50").Value = "  -15.94%  "

$cell = $ws.Range("D51")
$cell.NumberFormat = "@"
$cell.Value = "0.0209"
$cell.Style = "Normal"
$ws.Range("E51").Value = "  -13.92%  "
